$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Hunk 0: ALC row 18
$ws_ALC.Range("H18").Value = 679.1818
$ws_ALC.Range("I18").Value = 647.1
$ws_ALC.Range("J18").Value = 1000
$ws_ALC.Range("K18").Value = 647.1
$ws_ALC.Range("L18").Value = 1000
$ws_ALC.Range("M18").Value = -363.1
$ws_ALC.Range("N18").Value = -1568

# Hunk 1: ALC row 62
$ws_ALC.Range("H62").Value = 2772.5557
$ws_ALC.Range("I62").Value = 3399.889
$ws_ALC.Range("J62").Value = 2145.2222
$ws_ALC.Range("K62").Value = 3399.889
$ws_ALC.Range("L62").Value = 2145.2222
$ws_ALC.Range("M62").Value = -2775.889
$ws_ALC.Range("N62").Value = -3393.2222

# Hunk 2: ALC row 65
$ws_ALC.Range("H65").Value = 2772.5557
$ws_ALC.Range("I65").Value = 3399.889
$ws_ALC.Range("J65").Value = 2145.2222
$ws_ALC.Range("K65").Value = 16999.445
$ws_ALC.Range("L65").Value = 10726.111
$ws_ALC.Range("M65").Value = -13879.445
$ws_ALC.Range("N65").Value = -16966.111

# Hunk 3: ALC row 111
$ws_ALC.Range("H111").Value = 100723.5
$ws_ALC.Range("I111").Value = 642.6667
$ws_ALC.Range("J111").Value = 143615.28
$ws_ALC.Range("K111").Value = 1928.0001
$ws_ALC.Range("L111").Value = 430845.84
$ws_ALC.Range("M111").Value = 1138.9999
$ws_ALC.Range("N111").Value = -436979.84

# Hunk 4: ALC row 129
$ws_ALC.Range("H129").Value = 1114.9032
$ws_ALC.Range("I129").Value = 600
$ws_ALC.Range("J129").Value = 1141.0847
$ws_ALC.Range("K129").Value = 1800
$ws_ALC.Range("L129").Value = 3423.2541
$ws_ALC.Range("M129").Value = 3200
$ws_ALC.Range("N129").Value = -13423.2541

# Hunk 5: ALC row 132
$ws_ALC.Range("H132").Value = 1548.2778
$ws_ALC.Range("I132").Value = 1390.6522
$ws_ALC.Range("J132").Value = 2454.625
$ws_ALC.Range("K132").Value = 4171.9566
$ws_ALC.Range("L132").Value = 7363.875
$ws_ALC.Range("M132").Value = -1641.9566
$ws_ALC.Range("N132").Value = -12423.875

# Hunk 6: ALC row 137
$ws_ALC.Range("H137").Value = 1385.8182
$ws_ALC.Range("I137").Value = 1304.6364
$ws_ALC.Range("J137").Value = 1629.3636
$ws_ALC.Range("K137").Value = 3913.9092
$ws_ALC.Range("L137").Value = 4888.0908
$ws_ALC.Range("M137").Value = -1363.9092
$ws_ALC.Range("N137").Value = -9988.0908

# Hunk 7: ARM row 45
$ws_ARM.Range("H45").Value = 8436.615
$ws_ARM.Range("I45").Value = 10727.7
$ws_ARM.Range("J45").Value = 799.6667
$ws_ARM.Range("K45").Value = 10727.7
$ws_ARM.Range("L45").Value = 799.6667
$ws_ARM.Range("M45").Value = -10350.7
$ws_ARM.Range("N45").Value = -1553.6667

# Hunk 8: ARM row 74
$ws_ARM.Range("H74").Value = 1673.119
$ws_ARM.Range("I74").Value = 1330.3214
$ws_ARM.Range("K74").Value = 1330.3214
$ws_ARM.Range("M74").Value = -456.3214

# Hunk 9: ARM row 77
$ws_ARM.Range("H77").Value = 1673.119
$ws_ARM.Range("I77").Value = 1330.3214
$ws_ARM.Range("K77").Value = 6651.607
$ws_ARM.Range("M77").Value = -2283.607

# Hunk 10: ARM row 88
$ws_ARM.Range("H88").Value = 2350.3
$ws_ARM.Range("I88").Value = 2019.2
$ws_ARM.Range("K88").Value = 2019.2
$ws_ARM.Range("M88").Value = -1613.2

# Hunk 11: ARM row 91
$ws_ARM.Range("H91").Value = 2350.3
$ws_ARM.Range("I91").Value = 2019.2
$ws_ARM.Range("K91").Value = 2019.2
$ws_ARM.Range("M91").Value = -615.2

# Hunk 12: ARM row 109
$ws_ARM.Range("H109").Value = 44344.25
$ws_ARM.Range("J109").Value = 44344.25
$ws_ARM.Range("L109").Value = 44344.25
$ws_ARM.Range("N109").Value = -47118.25

# Hunk 13: ARM row 122
$ws_ARM.Range("H122").Value = 1117626.4
$ws_ARM.Range("I122").Value = 1223790.8
$ws_ARM.Range("J122").Value = 2900
$ws_ARM.Range("K122").Value = 3671372.4
$ws_ARM.Range("L122").Value = 8700
$ws_ARM.Range("M122").Value = -3668922.4
$ws_ARM.Range("N122").Value = -13600

# Hunk 14: ARM row 132
$ws_ARM.Range("H132").Value = 2328484.2
$ws_ARM.Range("I132").Value = 2022.4166
$ws_ARM.Range("J132").Value = 14293145
$ws_ARM.Range("K132").Value = 6067.2498
$ws_ARM.Range("L132").Value = 42879435
$ws_ARM.Range("M132").Value = -3537.2498
$ws_ARM.Range("N132").Value = -42884495

# Hunk 15: BSM row 20
$ws_BSM.Range("H20").Value = 8518.264999999999
$ws_BSM.Range("I20").Value = 1095.05
$ws_BSM.Range("J20").Value = 19122.857
$ws_BSM.Range("K20").Value = 1095.05
$ws_BSM.Range("L20").Value = 19122.857
$ws_BSM.Range("M20").Value = -848.05
$ws_BSM.Range("N20").Value = -19616.857

# Hunk 16: BSM row 94
$ws_BSM.Range("H94").Value = 1243.25
$ws_BSM.Range("I94").Value = 732.5454999999999
$ws_BSM.Range("J94").Value = 2366.8
$ws_BSM.Range("K94").Value = 732.5454999999999
$ws_BSM.Range("L94").Value = 2366.8
$ws_BSM.Range("M94").Value = -281.5454999999999
$ws_BSM.Range("N94").Value = -3268.8

# Hunk 17: BSM row 99
$ws_BSM.Range("H99").Value = 43479740
$ws_BSM.Range("I99").Value = 58824916
$ws_BSM.Range("K99").Value = 58824916
$ws_BSM.Range("M99").Value = -58823418

# Hunk 18: BSM row 107
$ws_BSM.Range("H107").Value = 1274.1666
$ws_BSM.Range("I107").Value = 1259
$ws_BSM.Range("J107").Value = 1350
$ws_BSM.Range("K107").Value = 1259
$ws_BSM.Range("L107").Value = 1350
$ws_BSM.Range("M107").Value = 661
$ws_BSM.Range("N107").Value = -5190

# Hunk 19: CRP row 22
$ws_CRP.Range("H22").Value = 154.28572
$ws_CRP.Range("I22").Value = 95
$ws_CRP.Range("J22").Value = 178
$ws_CRP.Range("K22").Value = 95
$ws_CRP.Range("L22").Value = 178
$ws_CRP.Range("M22").Value = 255
$ws_CRP.Range("N22").Value = -878

# Hunk 20: CRP row 134
$ws_CRP.Range("H134").Value = 259389.16
$ws_CRP.Range("I134").Value = 2960.2646
$ws_CRP.Range("J134").Value = 2003105.6
$ws_CRP.Range("K134").Value = 8880.793799999999
$ws_CRP.Range("L134").Value = 6009316.800000001
$ws_CRP.Range("M134").Value = -6345.793799999999
$ws_CRP.Range("N134").Value = -6014386.800000001

# Hunk 21: CUL row 35
$ws_CUL.Range("H35").Value = 0
$ws_CUL.Range("I35").Value = 0
$ws_CUL.Range("K35").Value = 0
$ws_CUL.Range("M35").ClearContents()

# Hunk 22: CUL row 115
$ws_CUL.Range("H115").Value = 1587.5
$ws_CUL.Range("I115").Value = 850
$ws_CUL.Range("J115").Value = 1833.3334
$ws_CUL.Range("K115").Value = 2550
$ws_CUL.Range("L115").Value = 5500.0002
$ws_CUL.Range("N115").Value = -7850.0002
$ws_CUL.Range("M115").Value = -1375

# Hunk 23: GSM row 2
$ws_GSM.Range("H2").Value = 242.75
$ws_GSM.Range("I2").Value = 208.16667
$ws_GSM.Range("J2").Value = 263.5
$ws_GSM.Range("K2").Value = 208.16667
$ws_GSM.Range("L2").Value = 263.5
$ws_GSM.Range("M2").Value = -95.16667000000001
$ws_GSM.Range("N2").Value = -489.5

# Hunk 24: GSM row 9
$ws_GSM.Range("H9").Value = 90500
$ws_GSM.Range("I9").Value = 90500
$ws_GSM.Range("K9").Value = 90500
$ws_GSM.Range("M9").Value = -90330

# Hunk 25: GSM row 80
$ws_GSM.Range("H80").Value = 6512.0835
$ws_GSM.Range("I80").Value = 9991.154
$ws_GSM.Range("J80").Value = 2400.4546
$ws_GSM.Range("K80").Value = 9991.154
$ws_GSM.Range("L80").Value = 2400.4546
$ws_GSM.Range("M80").Value = -8993.154
$ws_GSM.Range("N80").Value = -4396.4546

# Hunk 26: GSM row 83
$ws_GSM.Range("H83").Value = 6512.0835
$ws_GSM.Range("I83").Value = 9991.154
$ws_GSM.Range("J83").Value = 2400.4546
$ws_GSM.Range("K83").Value = 49955.77
$ws_GSM.Range("L83").Value = 12002.273
$ws_GSM.Range("M83").Value = -44963.77
$ws_GSM.Range("N83").Value = -21986.273

# Hunk 27: GSM row 122
$ws_GSM.Range("H122").Value = 52594876
$ws_GSM.Range("I122").Value = 56044980
$ws_GSM.Range("J122").Value = 41669548
$ws_GSM.Range("K122").Value = 168134940
$ws_GSM.Range("L122").Value = 125008644
$ws_GSM.Range("M122").Value = -168132490
$ws_GSM.Range("N122").Value = -125013544

# Hunk 28: LTW row 7
$ws_LTW.Range("H7").Value = 2663.125
$ws_LTW.Range("I7").Value = 2250
$ws_LTW.Range("J7").Value = 3902.5
$ws_LTW.Range("K7").Value = 2250
$ws_LTW.Range("L7").Value = 3902.5
$ws_LTW.Range("M7").Value = -2138
$ws_LTW.Range("N7").Value = -4126.5

# Hunk 29: LTW row 40
$ws_LTW.Range("H40").Value = 35716840
$ws_LTW.Range("I40").Value = 47621684
$ws_LTW.Range("J40").Value = 2314.8572
$ws_LTW.Range("K40").Value = 47621684
$ws_LTW.Range("L40").Value = 2314.8572
$ws_LTW.Range("M40").Value = -47621548
$ws_LTW.Range("N40").Value = -2586.8572

# Hunk 30: LTW row 61
$ws_LTW.Range("H61").Value = 1033.4231
$ws_LTW.Range("I61").Value = 962.76
$ws_LTW.Range("J61").Value = 2800
$ws_LTW.Range("K61").Value = 962.76
$ws_LTW.Range("L61").Value = 2800
$ws_LTW.Range("M61").Value = -760.76
$ws_LTW.Range("N61").Value = -3204

# Hunk 31: LTW row 82
$ws_LTW.Range("H82").Value = 531439.2
$ws_LTW.Range("I82").Value = 770963.0600000001
$ws_LTW.Range("J82").Value = 142212.88
$ws_LTW.Range("K82").Value = 770963.0600000001
$ws_LTW.Range("L82").Value = 142212.88
$ws_LTW.Range("M82").Value = -770602.0600000001
$ws_LTW.Range("N82").Value = -142934.88

# Hunk 32: LTW row 85
$ws_LTW.Range("H85").Value = 531439.2
$ws_LTW.Range("I85").Value = 770963.0600000001
$ws_LTW.Range("J85").Value = 142212.88
$ws_LTW.Range("K85").Value = 770963.0600000001
$ws_LTW.Range("L85").Value = 142212.88
$ws_LTW.Range("M85").Value = -769715.0600000001
$ws_LTW.Range("N85").Value = -144708.88

# Hunk 33: LTW row 113
$ws_LTW.Range("H113").Value = 1033.4231
$ws_LTW.Range("I113").Value = 962.76
$ws_LTW.Range("J113").Value = 2800
$ws_LTW.Range("K113").Value = 962.76
$ws_LTW.Range("L113").Value = 2800
$ws_LTW.Range("M113").Value = 1207.24
$ws_LTW.Range("N113").Value = -7140

# Hunk 34: LTW row 126
$ws_LTW.Range("H126").Value = 2663.125
$ws_LTW.Range("I126").Value = 2250
$ws_LTW.Range("J126").Value = 3902.5
$ws_LTW.Range("K126").Value = 6750
$ws_LTW.Range("L126").Value = 11707.5
$ws_LTW.Range("M126").Value = -4280
$ws_LTW.Range("N126").Value = -16647.5

# Hunk 35: WVR row 17
$ws_WVR.Range("H17").Value = 100
$ws_WVR.Range("I17").Value = 100
$ws_WVR.Range("K17").Value = 100
$ws_WVR.Range("M17").Value = 72

# Hunk 36: WVR row 81
$ws_WVR.Range("H81").Value = 1640.2667
$ws_WVR.Range("I81").Value = 1467
$ws_WVR.Range("J81").Value = 2333.3333
$ws_WVR.Range("K81").Value = 2934
$ws_WVR.Range("L81").Value = 4666.6666
$ws_WVR.Range("M81").Value = -1873
$ws_WVR.Range("N81").Value = -6788.6666

# Hunk 37: WVR row 84
$ws_WVR.Range("H84").Value = 1640.2667
$ws_WVR.Range("I84").Value = 1467
$ws_WVR.Range("J84").Value = 2333.3333
$ws_WVR.Range("K84").Value = 14670
$ws_WVR.Range("L84").Value = 23333.333
$ws_WVR.Range("M84").Value = -9366
$ws_WVR.Range("N84").Value = -33941.333

# Hunk 38: WVR row 122
$ws_WVR.Range("H122").Value = 1684.3125
$ws_WVR.Range("I122").Value = 1334.762
$ws_WVR.Range("J122").Value = 2351.6365
$ws_WVR.Range("K122").Value = 4004.286
$ws_WVR.Range("L122").Value = 7054.9095
$ws_WVR.Range("M122").Value = -1554.286
$ws_WVR.Range("N122").Value = -11954.9095

# Hunk 39: WVR row 126
$ws_WVR.Range("H126").Value = 813.7778
$ws_WVR.Range("I126").Value = 683.93335
$ws_WVR.Range("J126").Value = 1463
$ws_WVR.Range("K126").Value = 2051.80005
$ws_WVR.Range("L126").Value = 4389
$ws_WVR.Range("M126").Value = 418.1999500000002
$ws_WVR.Range("N126").Value = -9329
